$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "GM0011"
$ws.Range("B9").Value = "03/20/2025 06:25:21 PM"
$ws.Range("C9").Value = "ba kha 1111"
$ws.Range("D9").Value = 2100
$ws.Range("E9").Value = 9100
$ws.Range("F9").Value = 7000
